$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 319; existing rows 319.. shift down to 321..
$ws.Rows("319:320").Insert()

# --- New row 319: Early Majestic / Primera ---
$ws.Range("A319").Value = 5
$ws.Range("B319").Value = "Macroferia Regional de Talca"
$ws.Range("C319").Value = "Maule"
$ws.Range("D319").Value = 44889
$ws.Range("E319").Value = 7
$ws.Range("F319").Value = "Fruta"
$ws.Range("G319").Value = 100103
$ws.Range("H319").Value = "Frutos de hueso (carozo)"
$ws.Range("I319").Value = 100103004
$ws.Range("J319").Value = "Durazno"
$ws.Range("K319").Value = "Early Majestic"
$ws.Range("L319").Value = "Primera"
$ws.Range("M319").Value = 280
$ws.Range("N319").Value = 18000
$ws.Range("O319").Value = 18000
$ws.Range("P319").Value = 18000
$ws.Range("Q319").Value = "`$/caja 16 kilos granel"
$ws.Range("R319").Value = "Región de O'Higgins"
$ws.Range("S319").Value = 1125
$ws.Range("T319").Value = 16

# --- New row 320: Early Majestic / Segunda ---
$ws.Range("A320").Value = 5
$ws.Range("B320").Value = "Macroferia Regional de Talca"
$ws.Range("C320").Value = "Maule"
$ws.Range("D320").Value = 44889
$ws.Range("E320").Value = 7
$ws.Range("F320").Value = "Fruta"
$ws.Range("G320").Value = 100103
$ws.Range("H320").Value = "Frutos de hueso (carozo)"
$ws.Range("I320").Value = 100103004
$ws.Range("J320").Value = "Durazno"
$ws.Range("K320").Value = "Early Majestic"
$ws.Range("L320").Value = "Segunda"
$ws.Range("M320").Value = 350
$ws.Range("N320").Value = 10000
$ws.Range("O320").Value = 10000
$ws.Range("P320").Value = 10000
$ws.Range("Q320").Value = "`$/bandeja 8 kilos granel"
$ws.Range("R320").Value = "Región de O'Higgins"
$ws.Range("S320").Value = 1250
$ws.Range("T320").Value = 8
